$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 - heart_rhythm: Atrial Dysrhythmias
$ws.Range("C14").Value = 0
$ws.Range("D14").Value = "121 (21.1)"
$ws.Range("E14").Value = "127 (9.3)"
$ws.Range("F14").Value = "47 (7.4)"
$ws.Range("G14").Value = "431 (11.5)"
$ws.Range("H14").Value = "1886 (13.8)"

# Row 15 - Bradycardia with Pacer
$ws.Range("D15").Value = "39 (6.8)"
$ws.Range("E15").Value = "94 (6.9)"
$ws.Range("F15").Value = "62 (9.7)"
$ws.Range("G15").Value = "287 (7.6)"
$ws.Range("H15").Value = "1529 (11.2)"

# Row 16 - Bundle Branch Blocks
$ws.Range("D16").Value = "19 (3.3)"
$ws.Range("E16").Value = "52 (3.8)"
$ws.Range("F16").Value = "23 (3.6)"
$ws.Range("G16").Value = "156 (4.2)"
$ws.Range("H16").Value = "709 (5.2)"

# Row 17 - Normal (or near normal) Rhythm
$ws.Range("D17").Value = "394 (68.8)"
$ws.Range("E17").Value = "1092 (80.0)"
$ws.Range("F17").Value = "507 (79.3)"
$ws.Range("G17").Value = "2883 (76.7)"
$ws.Range("H17").Value = "9577 (69.9)"

# Row 18 - Ventricular Dysrhythmias (severe)
$ws.Range("H18").Value = "6 (0.0)"
